$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2282608695652174
$ws.Range("C2").Value = 0.4782608695652174
$ws.Range("J2").Value = 0.0108695652173913
$ws.Range("P2").Value = 0.1548913043478261
$ws.Range("S2").Value = 0.1277173913043478
$ws.Range("B3").Value = 0.0223463687150838
$ws.Range("C3").Value = 0.00558659217877095
$ws.Range("J3").Value = 0.09497206703910614
$ws.Range("O3").Value = 0.00558659217877095
$ws.Range("P3").Value = 0.6312849162011173
$ws.Range("S3").Value = 0.2402234636871508
$ws.Range("J4").Value = 0.07547169811320754
$ws.Range("P4").Value = 0.6037735849056604
$ws.Range("S4").Value = 0.3207547169811321
$ws.Range("B6").Value = 0.05797101449275362
$ws.Range("D6").Value = 0.00966183574879227
$ws.Range("F6").Value = 0.05797101449275362
$ws.Range("J6").Value = 0.2946859903381642
$ws.Range("O6").Value = 0.01932367149758454
$ws.Range("Q6").Value = 0.1352657004830918
$ws.Range("R6").Value = 0.0821256038647343
$ws.Range("S6").Value = 0.3429951690821256
$ws.Range("B7").Value = 0.1144067796610169
$ws.Range("D7").Value = 0.01694915254237288
$ws.Range("F7").Value = 0.03389830508474576
$ws.Range("J7").Value = 0.1610169491525424
$ws.Range("O7").Value = 0.0211864406779661
$ws.Range("Q7").Value = 0.173728813559322
$ws.Range("R7").Value = 0.07203389830508475
$ws.Range("S7").Value = 0.4067796610169492
$ws.Range("B8").Value = 0.1220043572984749
$ws.Range("D8").Value = 0.01525054466230937
$ws.Range("F8").Value = 0.06318082788671024
$ws.Range("J8").Value = 0.1263616557734205
$ws.Range("O8").Value = 0.0261437908496732
$ws.Range("Q8").Value = 0.1786492374727669
$ws.Range("R8").Value = 0.06753812636165578
$ws.Range("S8").Value = 0.4008714596949891
$ws.Range("B9").Value = 0.09770114942528736
$ws.Range("D9").Value = 0.02298850574712644
$ws.Range("F9").Value = 0.06896551724137931
$ws.Range("J9").Value = 0.1379310344827586
$ws.Range("O9").Value = 0.01149425287356322
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.09770114942528736
$ws.Range("S9").Value = 0.396551724137931
$ws.Range("B10").Value = 0.1204996326230713
$ws.Range("D10").Value = 0.02718589272593681
$ws.Range("E10").Value = 0.002939015429831006
$ws.Range("F10").Value = 0.05290227773695812
$ws.Range("J10").Value = 0.1102130786186628
$ws.Range("O10").Value = 0.01689933872152829
$ws.Range("Q10").Value = 0.2204261572373255
$ws.Range("R10").Value = 0.08229243203526819
$ws.Range("S10").Value = 0.3666421748714181
$ws.Range("G11").Value = 0.1550802139037433
$ws.Range("J11").Value = 0.1016042780748663
$ws.Range("K11").Value = 0.1978609625668449
$ws.Range("L11").Value = 0.5347593582887701
$ws.Range("S11").Value = 0.0106951871657754
$ws.Range("G12").Value = 0.7230046948356808
$ws.Range("J12").Value = 0.1643192488262911
$ws.Range("K12").Value = 0.01408450704225352
$ws.Range("L12").Value = 0.06572769953051644
$ws.Range("S12").Value = 0.03286384976525822
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.32
$ws.Range("S13").Value = 0.08
$ws.Range("F15").Value = 0.01298701298701299
$ws.Range("H15").Value = 0.2077922077922078
$ws.Range("I15").Value = 0.0303030303030303
$ws.Range("J15").Value = 0.3506493506493507
$ws.Range("K15").Value = 0.0735930735930736
$ws.Range("M15").Value = 0.01731601731601732
$ws.Range("O15").Value = 0.04761904761904762
$ws.Range("S15").Value = 0.2597402597402597
$ws.Range("F16").Value = 0.025
$ws.Range("H16").Value = 0.205
$ws.Range("I16").Value = 0.07000000000000001
$ws.Range("J16").Value = 0.36
$ws.Range("K16").Value = 0.17
$ws.Range("M16").Value = 0.01
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.12
$ws.Range("F17").Value = 0.02132196162046908
$ws.Range("H17").Value = 0.1663113006396588
$ws.Range("I17").Value = 0.06183368869936034
$ws.Range("J17").Value = 0.4243070362473347
$ws.Range("K17").Value = 0.11727078891258
$ws.Range("M17").Value = 0.0255863539445629
$ws.Range("N17").Value = 0.002132196162046908
$ws.Range("O17").Value = 0.08102345415778252
$ws.Range("S17").Value = 0.1002132196162047
$ws.Range("F18").Value = 0.03061224489795918
$ws.Range("H18").Value = 0.1836734693877551
$ws.Range("I18").Value = 0.1122448979591837
$ws.Range("J18").Value = 0.3928571428571428
$ws.Range("K18").Value = 0.1224489795918367
$ws.Range("M18").Value = 0.01020408163265306
$ws.Range("O18").Value = 0.06122448979591837
$ws.Range("S18").Value = 0.08673469387755102
$ws.Range("F19").Value = 0.01829268292682927
$ws.Range("H19").Value = 0.1890243902439024
$ws.Range("I19").Value = 0.08460365853658537
$ws.Range("J19").Value = 0.3879573170731707
$ws.Range("K19").Value = 0.1204268292682927
$ws.Range("M19").Value = 0.02362804878048781
$ws.Range("N19").Value = 0.0007621951219512195
$ws.Range("O19").Value = 0.06783536585365854
$ws.Range("S19").Value = 0.1074695121951219

Write-Output "Applied 110 cell updates"
